$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains exact text formatting (avoid numeric auto-conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '24.116.77'
$ws.Range('E2').Value = '  -3.34%  '
$ws.Range('D3').Value = '1.643.72'
$ws.Range('E3').Value = '  -3.35%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '308.22'
$ws.Range('E5').Value = '  -2.42%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = '0.3896'
$ws.Range('E7').Value = '  -1.79%  '
$ws.Range('D8').Value = '0.3870'
$ws.Range('E8').Value = '  -3.78%  '
$ws.Range('D9').Value = '1.003'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = '1.361'
$ws.Range('E10').Value = '  -7.17%  '
$ws.Range('D11').Value = '48.92'
$ws.Range('E11').Value = '  -7.32%  '
$ws.Range('D12').Value = '0.08469'
$ws.Range('E12').Value = '  -3.69%  '
$ws.Range('D13').Value = '24.20'
$ws.Range('E13').Value = '  -6.38%  '
$ws.Range('D14').Value = '7.159'
$ws.Range('E14').Value = '  -3.88%  '
$ws.Range('D15').Value = '0.00001286'
$ws.Range('E15').Value = '  -4.53%  '
$ws.Range('D16').Value = '7.519'
$ws.Range('E16').Value = '  -5.39%  '
$ws.Range('D17').Value = '1.653.17'
$ws.Range('E17').Value = '  -3.26%  '
$ws.Range('D18').Value = '94.78'
$ws.Range('E18').Value = '  -1.60%  '
$ws.Range('D19').Value = '0.06936'
$ws.Range('E19').Value = '  -3.58%  '
$ws.Range('D20').Value = '21.11'
$ws.Range('E20').Value = '  +2.59%  '
$ws.Range('D21').Value = '6.969'
$ws.Range('E21').Value = '  -5.05%  '
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('E23').Value = '  -3.69%  '
$ws.Range('D24').Value = '24.179.15'
$ws.Range('E24').Value = '  -3.07%  '
$ws.Range('D25').Value = '2.341'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '2.718'
$ws.Range('E26').Value = '  -7.13%  '
$ws.Range('D27').Value = '22.57'
$ws.Range('E27').Value = '  -4.85%  '
$ws.Range('D28').Value = '8.964'
$ws.Range('E28').Value = '  +8.01%  '
$ws.Range('D29').Value = '158.05'
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('D30').Value = '141.84'
$ws.Range('E30').Value = '  -5.23%  '
$ws.Range('D31').Value = '5.389'
$ws.Range('E31').Value = '  -12.41%  '
$ws.Range('D32').Value = '2.462'
$ws.Range('E32').Value = '  -6.72%  '
$ws.Range('D33').Value = '1.829.33'
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('D34').Value = '7.187'
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('D35').Value = '0.08072'
$ws.Range('E35').Value = '  -5.48%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.9894'
$ws.Range('E36').Value = '  -4.76%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02952'
$ws.Range('E37').Value = '  -5.67%  '
$ws.Range('D38').Value = '0.2716'
$ws.Range('E38').Value = '  -4.62%  '
$ws.Range('D39').Value = '0.09291'
$ws.Range('E39').Value = '  -2.47%  '
$ws.Range('D40').Value = '1.481'
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('D41').Value = '10.06'
$ws.Range('E41').Value = '  -7.39%  '
$ws.Range('D42').Value = '0.7651'
$ws.Range('E42').Value = '  -6.81%  '
$ws.Range('D43').Value = '13.08'
$ws.Range('E43').Value = '  -6.01%  '
$ws.Range('D44').Value = '16.07'
$ws.Range('E44').Value = '  -6.02%  '
$ws.Range('D45').Value = '2.495'
$ws.Range('E45').Value = '  -6.50%  '
$ws.Range('D46').Value = '0.6889'
$ws.Range('E46').Value = '  -6.52%  '
$ws.Range('D47').Value = '4.101'
$ws.Range('E47').Value = '  -3.50%  '
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = '0.08416'
$ws.Range('E49').Value = '  -3.70%  '
$ws.Range('D50').Value = '134.29'
$ws.Range('E50').Value = '  -3.36%  '
$ws.Range('D51').Value = '1.263'
$ws.Range('E51').Value = '  -9.79%  '
